# Apply the "3e version avec organisation fichiers" update to
# StructureDefinition-SituationOperationnelle.xlsx
#
#  - bump the generation Date on the Metadata sheet
#  - point the three ValueSet bindings at the new terminology-server
#    CodeSystem URLs instead of the old mos/ValueSet ones
#  - drop the "ContactPoint { ... }" wrapper around the Telecommunication
#    structure-definition URL
#  - column widths on the Elements sheet shift slightly because the cell
#    contents above changed length (best-fit autosize)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
$meta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Elements sheet ---------------------------------------------------
# SituationOperationnelle.modeExerciceOffre -> Binding Value Set
$elements.Range("Z4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R23-ModeExercice?vs"

# SituationOperationnelle.competenceSpecifique -> Binding Value Set
$elements.Range("Z5").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R243-CompetenceSpecifique?vs"

# SituationOperationnelle.secteurConventionnement -> Binding Value Set
$elements.Range("Z7").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R282-CNAMAmeliSecteurConventionnement?vs"

# SituationOperationnelle.telecommunication -> Type(s)
$elements.Range("K10").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/Telecommunication`n"

# Column best-fit widths follow the content changes above.
$elements.Columns.Item(11).ColumnWidth = 63.6666666666667
$elements.Columns.Item(26).ColumnWidth = 86.6666666666667
